$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36; this shifts existing rows 36..78 down to 37..79
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new record
$ws.Cells.Item(36, 1).Value = 9
$ws.Cells.Item(36, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(36, 3).Value = "Metropolitana"
$ws.Cells.Item(36, 4).Value = 44818
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(36, 6).Value = 100112029
$ws.Cells.Item(36, 7).Value = "Orégano"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 16
$ws.Cells.Item(36, 11).Value = 18000
$ws.Cells.Item(36, 12).Value = 18000
$ws.Cells.Item(36, 13).Value = 18000
$ws.Cells.Item(36, 14).Value = "$/docena de atados"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 6000
$ws.Cells.Item(36, 17).Value = 3
$ws.Cells.Item(36, 18).Value = "Hortaliza"
